$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 38000
$ws.Range("J93").Value = 38000
$ws.Range("L93").Value = 38000
$ws.Range("N93").Value = -42992
$ws.Range("H132").Value = 4078.2917
$ws.Range("I132").Value = 4380.864
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 13142.592
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = -10612.592
$ws.Range("N132").Value = -7310
$ws.Range("H138").Value = 3377.01
$ws.Range("I138").Value = 2592.2666
$ws.Range("J138").Value = 3713.3286
$ws.Range("K138").Value = 7776.7998
$ws.Range("L138").Value = 11139.9858
$ws.Range("M138").Value = -2636.7998
$ws.Range("N138").Value = -21419.9858

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1471851
$ws.Range("I2").Value = 1474.5714
$ws.Range("J2").Value = 4902729.5
$ws.Range("K2").Value = 1474.5714
$ws.Range("L2").Value = 4902729.5
$ws.Range("M2").Value = -1361.5714
$ws.Range("N2").Value = -4902955.5
$ws.Range("H5").Value = 266.25
$ws.Range("I5").Value = 213.75
$ws.Range("K5").Value = 213.75
$ws.Range("M5").Value = -101.75
$ws.Range("H32").Value = 20716.234
$ws.Range("I32").Value = 15721.362
$ws.Range("J32").Value = 69000
$ws.Range("K32").Value = 15721.362
$ws.Range("L32").Value = 69000
$ws.Range("M32").Value = -15434.362
$ws.Range("N32").Value = -69574
$ws.Range("H80").Value = 20983.334
$ws.Range("J80").Value = 20983.334
$ws.Range("L80").Value = 20983.334
$ws.Range("N80").Value = -22979.334
$ws.Range("H83").Value = 20983.334
$ws.Range("J83").Value = 20983.334
$ws.Range("L83").Value = 62950.00199999999
$ws.Range("N83").Value = -72934.00199999999
$ws.Range("H97").Value = 899.80646
$ws.Range("I97").Value = 846.3077
$ws.Range("J97").Value = 1178
$ws.Range("K97").Value = 846.3077
$ws.Range("L97").Value = 1178
$ws.Range("M97").Value = -350.3077
$ws.Range("N97").Value = -2170
$ws.Range("H102").Value = 2418.3333
$ws.Range("I102").Value = 2352.5
$ws.Range("K102").Value = 2352.5
$ws.Range("M102").Value = -730.5
$ws.Range("H111").Value = 48429.332
$ws.Range("J111").Value = 48429.332
$ws.Range("L111").Value = 48429.332
$ws.Range("N111").Value = -56609.332
$ws.Range("H116").Value = 1471851
$ws.Range("I116").Value = 1474.5714
$ws.Range("J116").Value = 4902729.5
$ws.Range("K116").Value = 1474.5714
$ws.Range("L116").Value = 4902729.5
$ws.Range("M116").Value = 819.4286
$ws.Range("N116").Value = -4907317.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1471851
$ws.Range("I3").Value = 1474.5714
$ws.Range("J3").Value = 4902729.5
$ws.Range("K3").Value = 1474.5714
$ws.Range("L3").Value = 4902729.5
$ws.Range("M3").Value = -1360.5714
$ws.Range("N3").Value = -4902957.5
$ws.Range("H4").Value = 266.25
$ws.Range("I4").Value = 213.75
$ws.Range("K4").Value = 213.75
$ws.Range("M4").Value = -98.75
$ws.Range("H15").Value = 12500
$ws.Range("J15").Value = 12500
$ws.Range("L15").Value = 12500
$ws.Range("N15").Value = -12954
$ws.Range("H86").Value = 1563.5238
$ws.Range("I86").Value = 1402.1875
$ws.Range("J86").Value = 2079.8
$ws.Range("K86").Value = 1402.1875
$ws.Range("L86").Value = 2079.8
$ws.Range("M86").Value = -279.1875
$ws.Range("N86").Value = -4325.8
$ws.Range("H89").Value = 1563.5238
$ws.Range("I89").Value = 1402.1875
$ws.Range("J89").Value = 2079.8
$ws.Range("K89").Value = 7010.9375
$ws.Range("L89").Value = 10399
$ws.Range("M89").Value = -1394.9375
$ws.Range("N89").Value = -21631
$ws.Range("H94").Value = 245.6
$ws.Range("I94").Value = 245.6
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 245.6
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 205.4
$ws.Range("N94").ClearContents()
$ws.Range("H141").Value = 61827.668
$ws.Range("J141").Value = 57391.43
$ws.Range("L141").Value = 57391.43
$ws.Range("N141").Value = -67751.42999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17232.605
$ws.Range("I31").Value = 35173.8
$ws.Range("J31").Value = 7621.25
$ws.Range("K31").Value = 35173.8
$ws.Range("L31").Value = 7621.25
$ws.Range("M31").Value = -34878.8
$ws.Range("N31").Value = -8211.25
$ws.Range("H34").Value = 17232.605
$ws.Range("I34").Value = 35173.8
$ws.Range("J34").Value = 7621.25
$ws.Range("K34").Value = 35173.8
$ws.Range("L34").Value = 7621.25
$ws.Range("M34").Value = -34971.8
$ws.Range("N34").Value = -8025.25
$ws.Range("H51").Value = 9999
$ws.Range("J51").Value = 10298.75
$ws.Range("L51").Value = 10298.75
$ws.Range("N51").Value = -11770.75
$ws.Range("H61").Value = 9999
$ws.Range("J61").Value = 10298.75
$ws.Range("L61").Value = 10298.75
$ws.Range("N61").Value = -10994.75
$ws.Range("H68").Value = 18147.545
$ws.Range("J68").Value = 18535.5
$ws.Range("L68").Value = 18535.5
$ws.Range("N68").Value = -20033.5
$ws.Range("H71").Value = 18147.545
$ws.Range("J71").Value = 18535.5
$ws.Range("L71").Value = 55606.5
$ws.Range("N71").Value = -63094.5
$ws.Range("H133").Value = 34975.332
$ws.Range("J133").Value = 34975.332
$ws.Range("L133").Value = 34975.332
$ws.Range("N133").Value = -40035.332
$ws.Range("H141").Value = 44158.8
$ws.Range("I141").Value = 17098.666
$ws.Range("J141").Value = 50923.832
$ws.Range("K141").Value = 17098.666
$ws.Range("L141").Value = 50923.832
$ws.Range("M141").Value = -11918.666
$ws.Range("N141").Value = -61283.832

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6158843
$ws.Range("I122").Value = 12821088
$ws.Range("J122").Value = 745769
$ws.Range("K122").Value = 115389792
$ws.Range("L122").Value = 6711921
$ws.Range("M122").Value = -115387342
$ws.Range("N122").Value = -6716821
$ws.Range("H131").Value = 46671690
$ws.Range("I131").Value = 111121320
$ws.Range("K131").Value = 333363960
$ws.Range("M131").Value = -333358920
$ws.Range("H133").Value = 5700
$ws.Range("I133").Value = 5375
$ws.Range("J133").Value = 7000
$ws.Range("K133").Value = 16125
$ws.Range("L133").Value = 21000
$ws.Range("M133").Value = -11065
$ws.Range("N133").Value = -31120
$ws.Range("H139").Value = 1752.4706
$ws.Range("I139").Value = 1797.4375
$ws.Range("J139").Value = 1033
$ws.Range("K139").Value = 5392.3125
$ws.Range("L139").Value = 3099
$ws.Range("M139").Value = -252.3125
$ws.Range("N139").Value = -13379

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1842
$ws.Range("I97").Value = 1768.8889
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 1768.8889
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -1272.8889
$ws.Range("N97").Value = -3492
$ws.Range("H132").Value = 1694.037
$ws.Range("I132").Value = 1226.35
$ws.Range("K132").Value = 3679.05
$ws.Range("M132").Value = -1149.05

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17301
$ws.Range("I132").Value = 30576
$ws.Range("J132").Value = 4026
$ws.Range("K132").Value = 91728
$ws.Range("L132").Value = 12078
$ws.Range("M132").Value = -89198
$ws.Range("N132").Value = -17138
